$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.931.67"
$ws.Range("E2").Value = "  -1.97%  "
$ws.Range("D3").Value = "3.098.50"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "525.90"
$ws.Range("E5").Value = "  +0.53%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.31"
$ws.Range("E6").Value = "  -1.94%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "3.096.96"
$ws.Range("E8").Value = "  -0.31%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.444"
$ws.Range("E10").Value = "  -3.08%  "
$ws.Range("E11").Value = "  -1.29%  "
$ws.Range("E12").Value = "  +1.88%  "
$ws.Range("D13").Value = "3.631.43"
$ws.Range("E13").Value = "  -0.26%  "
$ws.Range("E14").Value = "  +2.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.63"
$ws.Range("E15").Value = "  -5.33%  "
$ws.Range("D17").Value = "57.988.55"
$ws.Range("E17").Value = "  -1.79%  "
$ws.Range("D18").Value = "3.101.84"
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.10"
$ws.Range("E19").Value = "  -1.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.69"
$ws.Range("E20").Value = "  -2.61%  "
$ws.Range("E21").Value = "  -2.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "342.70"
$ws.Range("E22").Value = "  -0.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("E24").Value = "  +0.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "67.50"
$ws.Range("E25").Value = "  +2.36%  "
$ws.Range("E26").Value = "  -0.34%  "
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("D28").Value = "0.0₃0925"
$ws.Range("E28").Value = "  -1.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.36"
$ws.Range("E30").Value = "  -6.39%  "
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("E32").Value = "  +1.10%  "
$ws.Range("E33").Value = "  -0.57%  "
$ws.Range("E34").Value = "  -2.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "159.12"
$ws.Range("E35").Value = "  +2.61%  "
$ws.Range("E36").Value = "  -0.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.15"
$ws.Range("E37").Value = "  -0.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.15"
$ws.Range("E38").Value = "  -3.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.25"
$ws.Range("E39").Value = "  -4.73%  "
$ws.Range("E40").Value = "  -2.80%  "
$ws.Range("E41").Value = "  +7.95%  "
$ws.Range("E42").Value = "  +1.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.685"
$ws.Range("E43").Value = "  +2.95%  "
$ws.Range("D44").Value = "3.139.36"
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "36.85"
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0262"
$ws.Range("E47").Value = "  +1.45%  "
$ws.Range("D48").Value = "2.267.76"
$ws.Range("E48").Value = "  -0.98%  "
$ws.Range("E49").Value = "  +2.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.11"
$ws.Range("E50").Value = "  +1.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.45"
$ws.Range("E51").Value = "  -2.26%  "
